$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (LinearRegression) - update C2 and D2
$ws.Range("C2").Value = 5504813121739529
$ws.Range("D2").Value = 5504813121739529

# Row 3 (RandomForestRegressor) - values updated
$ws.Range("B3").Value = 4740948609203.721
$ws.Range("C3").Value = 23747053742538.62
$ws.Range("D3").Value = 19268821053017.39

# Row 4 - model name changed GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 4391195768686.856
$ws.Range("C4").Value = 4626125326713.485
$ws.Range("D4").Value = 4625735725456.395

# Row 5 - model name changed AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 180292811073679.2
$ws.Range("C5").Value = 182141081921002.4
$ws.Range("D5").Value = 253657987019571
